$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current (before) layout, columns A..F:
#   A=kitchens_1  B=kitchens_2  C=living_rooms_1  D=bedrooms_1  E=living_rooms_2  F=bedrooms_2
# Target (after) layout, columns A..F:
#   A=living_rooms_1  B=bedrooms_1  C=kitchens_1  D=living_rooms_2  E=bedrooms_2  F=kitchens_2
#
# This is a pure re-ordering of the 6 columns (header + the 6 data rows move together).
# New column i (1-based) takes its full content from old column sourceCols[i-1].
$sourceCols = @(3, 4, 1, 5, 6, 2)

$lastRow = 7
$lastCol = 6

# Snapshot all original values (including header row) before overwriting anything.
$original = @{}
for ($c = 1; $c -le $lastCol; $c++) {
    $colValues = @()
    for ($r = 1; $r -le $lastRow; $r++) {
        $colValues += ,$ws.Cells.Item($r, $c).Value()
    }
    $original[$c] = $colValues
}

for ($newCol = 1; $newCol -le $lastCol; $newCol++) {
    $oldCol = $sourceCols[$newCol - 1]
    $colValues = $original[$oldCol]
    for ($r = 1; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, $newCol).Value = $colValues[$r - 1]
    }
}
